$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 12 with the new "admin screen" protocol entry
$ws.Range("A12").Value = "Serveur"
$ws.Range("B12").Value = "Client"
$ws.Range("C12").Value = "AD"
$ws.Range("E12").Value = "Envoi le client vers l'écran du mode admin"

# Update selection to F12
$ws.Range("F12").Select()
